# Support local pgs catalogs
# Update the sample "Scores" sheet so the example row points at a locally
# hosted PGS catalog entry instead of an external FTP download link.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2 values -------------------------------------------------
# A2 ("Polygenic Score (PGS) ID") stays "IBK000001" - no change needed.
# B2 ("PGS Name") : "Wu et al. - 2021 " -> "Copy of PGS000001"
$ws.Range("B2").Value = "Copy of PGS000001"

# D2 ("Number of Variants") : 5114 -> 77
$ws.Range("D2").Value = 77

# E2 ("FTP link") : drop the external hyperlink, replace with a local file name
$ws.Hyperlinks.Delete()
$ws.Range("E2").Value = "IBK000001.txt.gz"
$ws.Range("E2").Style = "Normal"

# Remove the now-unused "Hyperlink" cell style so it no longer lingers in the
# workbook's style table.
$wb.Styles.Item("Hyperlink").Delete()

# --- Restore the active selection ----------------------------------------
$ws.Range("D3").Select()
